# Add 2022-Q3 data
# 1) Insert a new worksheet "2022-Q3" right after "总计" and before "2022-Q1",
#    populated with the fund-holding detail for the quarter.
# 2) Insert a new row into "总计" (the summary sheet) for "2022-Q3",
#    shifting the existing quarter rows down by one.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)
$q1 = $wb.Worksheets.Item("2022-Q1")

# --- Step 1: create the new "2022-Q3" sheet -------------------------------
# Duplicate the "2022-Q1" sheet (so the new sheet inherits the same column
# headers/styles/borders) right after "总计", then trim it down to a single
# data row and overwrite the values.
$q1.Copy($null, $total)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The copied sheet has 3 data rows (rows 2-4); keep only row 2 and remove
# the rest so the sheet ends up with exactly one holding.
$q3.Rows.Item(3).Delete()
$q3.Rows.Item(3).Delete()

$q3.Range("A2").Value = 0

$q3.Range("B2").NumberFormat = "@"
$q3.Range("B2").Value = "002123"

$q3.Range("C2").Value = "北信瑞丰外延增长主题灵活配置混合"

$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "0.15"

$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "88.55"

$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "3.95"

$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.0059"

$q3.Range("H2").Value = 8

# --- Step 2: add the 2022-Q3 row to the "总计" (summary) sheet ------------
$total.Rows.Item(2).Insert()

# New row comes in unformatted; clear it and copy the row-format from the
# row below (which holds the same per-row styling as every other data row).
$total.Range("A2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01

# Column A is just a running 0-based row index; refresh it for every row
# now that the table has grown by one entry.
$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
